# Apply odds update to Jogos da Semana FlashScore workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 3
$ws.Range("I2").Value = 2.55
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 1.91
$ws.Range("L2").Value = 3.4
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 1.95
$ws.Range("Y2").Value = 1.54
$ws.Range("AC2").Value = 7.5
$ws.Range("AD2").Value = 13
$ws.Range("AJ2").Value = 6
$ws.Range("AP2").Value = 11
$ws.Range("AQ2").Value = 26

# Row 5 update
$ws.Range("AM5").Value = 900
